$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201, shifting existing rows 201-253 down to 202-254
$ws.Rows.Item(201).Insert()

# Populate the new row 201 with the new weekly data point
$ws.Cells.Item(201, 1).Value = 3
$ws.Cells.Item(201, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(201, 3).Value = "Coquimbo"
$ws.Cells.Item(201, 4).Value = 44551
$ws.Cells.Item(201, 4).NumberFormat = $ws.Cells.Item(202, 4).NumberFormat
$ws.Cells.Item(201, 5).Value = 5
$ws.Cells.Item(201, 6).Value = 100112009
$ws.Cells.Item(201, 7).Value = "Acelga"
$ws.Cells.Item(201, 8).Value = "Sin especificar"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 280
$ws.Cells.Item(201, 11).Value = 2400
$ws.Cells.Item(201, 12).Value = 2500
$ws.Cells.Item(201, 13).Value = 2457
$ws.Cells.Item(201, 14).Value = "`$/docena de atados (6 kilos)"
$ws.Cells.Item(201, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(201, 16).Value = 410
$ws.Cells.Item(201, 17).Value = 6
$ws.Cells.Item(201, 18).Value = "Hortaliza"
